$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ------------------------------------------------------------------
# 1) Remove the stray "#truck" sheet (its data lives on "trucks" too)
# ------------------------------------------------------------------
$wsTruck = $wb.Worksheets.Item("#truck")
$wsTruck.Delete()

# ------------------------------------------------------------------
# 2) "car" sheet: add the new Toyota Camry row
# ------------------------------------------------------------------
$wsCar = $wb.Worksheets.Item("car")
$wsCar.Range("A4").Value = 3456778
$wsCar.Range("B4").Value = "Toyota"
$wsCar.Range("C4").Value = "Camry"
$wsCar.Range("D4").Value = 230
$wsCar.Range("E4").Value = 34000.8
$wsCar.Range("F4").Value = 2
$wsCar.Range("G4").Value = "COUPE"

# ------------------------------------------------------------------
# 3) "trucks" sheet: add the new DAF XF-105 row, matching the
#    number-format already used on the "Koegel"/"BPW" trailer cells
# ------------------------------------------------------------------
$wsTrucks = $wb.Worksheets.Item("trucks")
$wsTrucks.Range("F2:G2").NumberFormat = "@"

$wsTrucks.Range("A3").Value = 343890
$wsTrucks.Range("B3").Value = "DAF"
$wsTrucks.Range("C3").Value = "XF-105"
$wsTrucks.Range("D3").Value = 320
$wsTrucks.Range("E3").Value = 65056.6
$wsTrucks.Range("F3").Value = "Schmitz"
$wsTrucks.Range("G3").Value = "B1997"
$wsTrucks.Range("H3").Value = 20
$wsTrucks.Range("I3").Value = 92
$wsTrucks.Range("J3").Value = "TILTCOVEREDTRUCK"
$wsTrucks.Range("K3").Value = "TOP"
$wsTrucks.Range("F3:G3").NumberFormat = "@"

# ------------------------------------------------------------------
# 4) Make "trucks" the active sheet/tab and restore the selections
#    that were left behind after the data entry
# ------------------------------------------------------------------
$wsCar.Range("E5").Select()
$wsTrucks.Activate()
$wsTrucks.Range("G4").Select()
